$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (keep header row 1) so the shared-string table
# can be rebuilt from scratch in the exact order required by the target file.
$ws.Range("A2:T7").ClearContents()

# Fill the string columns column-by-column (A, then B, then C, then D) across
# all 12 data rows so new shared strings get interned in the needed order:
# M1, M2, Ifnk, Ifnar1, ECs, FAPs, Neutro, sCs

# Column A
$ws.Range("A2").Value = "M1"
$ws.Range("A3").Value = "M1"
$ws.Range("A4").Value = "M1"
$ws.Range("A5").Value = "M1"
$ws.Range("A6").Value = "M1"
$ws.Range("A7").Value = "M1"
$ws.Range("A8").Value = "M2"
$ws.Range("A9").Value = "M2"
$ws.Range("A10").Value = "M2"
$ws.Range("A11").Value = "M2"
$ws.Range("A12").Value = "M2"
$ws.Range("A13").Value = "M2"

# Column B
$ws.Range("B2").Value = "Ifnk"
$ws.Range("B3").Value = "Ifnk"
$ws.Range("B4").Value = "Ifnk"
$ws.Range("B5").Value = "Ifnk"
$ws.Range("B6").Value = "Ifnk"
$ws.Range("B7").Value = "Ifnk"
$ws.Range("B8").Value = "Ifnk"
$ws.Range("B9").Value = "Ifnk"
$ws.Range("B10").Value = "Ifnk"
$ws.Range("B11").Value = "Ifnk"
$ws.Range("B12").Value = "Ifnk"
$ws.Range("B13").Value = "Ifnk"

# Column C
$ws.Range("C2").Value = "Ifnar1"
$ws.Range("C3").Value = "Ifnar1"
$ws.Range("C4").Value = "Ifnar1"
$ws.Range("C5").Value = "Ifnar1"
$ws.Range("C6").Value = "Ifnar1"
$ws.Range("C7").Value = "Ifnar1"
$ws.Range("C8").Value = "Ifnar1"
$ws.Range("C9").Value = "Ifnar1"
$ws.Range("C10").Value = "Ifnar1"
$ws.Range("C11").Value = "Ifnar1"
$ws.Range("C12").Value = "Ifnar1"
$ws.Range("C13").Value = "Ifnar1"

# Column D
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "M1"
$ws.Range("D5").Value = "M2"
$ws.Range("D6").Value = "Neutro"
$ws.Range("D7").Value = "sCs"
$ws.Range("D8").Value = "ECs"
$ws.Range("D9").Value = "FAPs"
$ws.Range("D10").Value = "M1"
$ws.Range("D11").Value = "M2"
$ws.Range("D12").Value = "Neutro"
$ws.Range("D13").Value = "sCs"

# Fill numeric columns row-by-row

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1379033333333333
$ws.Range("H2").Value = 0.41371
$ws.Range("I2").Value = 0.1151434304008604
$ws.Range("J2").Value = 0.1151434304008603
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.568057
$ws.Range("N2").Value = 53.136114
$ws.Range("O2").Value = 0.1769362593421452
$ws.Range("P2").Value = 0.1285250503302229
$ws.Range("Q2").Value = 3.663823620490001
$ws.Range("R2").Value = 21.98294172294
$ws.Range("S2").Value = 0.02037304786295088
$ws.Range("T2").Value = 0.01479881518746509

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1379033333333333
$ws.Range("H3").Value = 0.41371
$ws.Range("I3").Value = 0.1151434304008604
$ws.Range("J3").Value = 0.1151434304008603
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 13.56641266666666
$ws.Range("N3").Value = 40.69923799999999
$ws.Range("O3").Value = 0.09034873381715106
$ws.Range("P3").Value = 0.09844287093240803
$ws.Range("Q3").Value = 1.870853528108889
$ws.Range("R3").Value = 16.83768175298
$ws.Range("S3").Value = 0.01040306314408099
$ws.Range("T3").Value = 0.0113350498576666

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1379033333333333
$ws.Range("H4").Value = 0.41371
$ws.Range("I4").Value = 0.1151434304008604
$ws.Range("J4").Value = 0.1151434304008603
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 39.790493
$ws.Range("N4").Value = 119.371479
$ws.Range("O4").Value = 0.2649941991919023
$ws.Range("P4").Value = 0.2887344254506106
$ws.Range("Q4").Value = 5.487241619676667
$ws.Range("R4").Value = 49.38517457709
$ws.Range("S4").Value = 0.03051234113128453
$ws.Range("T4").Value = 0.03324587222120478

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1379033333333333
$ws.Range("H5").Value = 0.41371
$ws.Range("I5").Value = 0.1151434304008604
$ws.Range("J5").Value = 0.1151434304008603
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 41.751761
$ws.Range("N5").Value = 125.255283
$ws.Range("O5").Value = 0.2780557273076938
$ws.Range("P5").Value = 0.3029661060969064
$ws.Range("Q5").Value = 5.757707014436668
$ws.Range("R5").Value = 51.81936312993
$ws.Range("S5").Value = 0.03201629028481405
$ws.Range("T5").Value = 0.03488455675118881

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1379033333333333
$ws.Range("H6").Value = 0.41371
$ws.Range("I6").Value = 0.1151434304008604
$ws.Range("J6").Value = 0.1151434304008603
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 18.009171
$ws.Range("N6").Value = 54.02751300000001
$ws.Range("O6").Value = 0.1199363337180826
$ws.Range("P6").Value = 0.1306811564643544
$ws.Range("Q6").Value = 2.483524711470001
$ws.Range("R6").Value = 22.35172240323
$ws.Range("S6").Value = 0.01380988089400241
$ws.Range("T6").Value = 0.01504707664405733

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1379033333333333
$ws.Range("H7").Value = 0.41371
$ws.Range("I7").Value = 0.1151434304008604
$ws.Range("J7").Value = 0.1151434304008603
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.470196
$ws.Range("N7").Value = 20.940392
$ws.Range("O7").Value = 0.06972874662302521
$ws.Range("P7").Value = 0.05065039072549785
$ws.Range("Q7").Value = 1.443874929053333
$ws.Range("R7").Value = 8.66324957432
$ws.Range("S7").Value = 0.00802880708372753
$ws.Range("T7").Value = 0.005832059739277745

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.059762333333333
$ws.Range("H8").Value = 3.179287
$ws.Range("I8").Value = 0.8848565695991397
$ws.Range("J8").Value = 0.8848565695991396
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 26.568057
$ws.Range("N8").Value = 53.136114
$ws.Range("O8").Value = 0.1769362593421452
$ws.Range("P8").Value = 0.1285250503302229
$ws.Range("Q8").Value = 28.155826078453
$ws.Range("R8").Value = 168.934956470718
$ws.Range("S8").Value = 0.1565632114791944
$ws.Range("T8").Value = 0.1137262351427578

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.059762333333333
$ws.Range("H9").Value = 3.179287
$ws.Range("I9").Value = 0.8848565695991397
$ws.Range("J9").Value = 0.8848565695991396
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 13.56641266666666
$ws.Range("N9").Value = 40.69923799999999
$ws.Range("O9").Value = 0.09034873381715106
$ws.Range("P9").Value = 0.09844287093240803
$ws.Range("Q9").Value = 14.37717314258955
$ws.Range("R9").Value = 129.394558283306
$ws.Range("S9").Value = 0.07994567067307007
$ws.Range("T9").Value = 0.08710782107474142

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.059762333333333
$ws.Range("H10").Value = 3.179287
$ws.Range("I10").Value = 0.8848565695991397
$ws.Range("J10").Value = 0.8848565695991396
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 39.790493
$ws.Range("N10").Value = 119.371479
$ws.Range("O10").Value = 0.2649941991919023
$ws.Range("P10").Value = 0.2887344254506106
$ws.Range("Q10").Value = 42.16846570616366
$ws.Range("R10").Value = 379.516191355473
$ws.Range("S10").Value = 0.2344818580606177
$ws.Range("T10").Value = 0.2554885532294058

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.059762333333333
$ws.Range("H11").Value = 3.179287
$ws.Range("I11").Value = 0.8848565695991397
$ws.Range("J11").Value = 0.8848565695991396
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 41.751761
$ws.Range("N11").Value = 125.255283
$ws.Range("O11").Value = 0.2780557273076938
$ws.Range("P11").Value = 0.3029661060969064
$ws.Range("Q11").Value = 44.24694365813566
$ws.Range("R11").Value = 398.222492923221
$ws.Range("S11").Value = 0.2460394370228798
$ws.Range("T11").Value = 0.2680815493457175

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.059762333333333
$ws.Range("H12").Value = 3.179287
$ws.Range("I12").Value = 0.8848565695991397
$ws.Range("J12").Value = 0.8848565695991396
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 18.009171
$ws.Range("N12").Value = 54.02751300000001
$ws.Range("O12").Value = 0.1199363337180826
$ws.Range("P12").Value = 0.1306811564643544
$ws.Range("Q12").Value = 19.085441080359
$ws.Range("R12").Value = 171.768969723231
$ws.Range("S12").Value = 0.1061264528240802
$ws.Range("T12").Value = 0.1156340798202971

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.059762333333333
$ws.Range("H13").Value = 3.179287
$ws.Range("I13").Value = 0.8848565695991397
$ws.Range("J13").Value = 0.8848565695991396
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 10.470196
$ws.Range("N13").Value = 20.940392
$ws.Range("O13").Value = 0.06972874662302521
$ws.Range("P13").Value = 0.05065039072549785
$ws.Range("Q13").Value = 11.09591934341733
$ws.Range("R13").Value = 66.575516060504
$ws.Range("S13").Value = 0.06169993953929769
$ws.Range("T13").Value = 0.04481833098622011
